# "Docs: 更新excel檔" (re-applied) — insert a new example row before the
# existing row 6 ("賣掉/nv/100"), containing the 虛詞/被/passive entry,
# highlighted with a red box border + red font, and bump the sheet's
# default row spacing slightly (15 -> 15.75), matching the restored
# worksheet snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the sheet-wide row height first (rows 1-4 and 7-10 end up at the
# new default 15.75 once the insert below shifts everything down).
$ws.UsedRange.RowHeight = 15.75

# Push rows 6-9 down to 7-10, opening up a blank row 6 for the new entry.
$ws.Rows.Item(6).Insert()

# Populate the new row.
$a6 = $ws.Range("A6")
$a6.Value = "虛詞"
$b6 = $ws.Range("B6")
$b6.Value = "被"
$c6 = $ws.Range("C6")
$c6.Value = "passive"

# Slightly taller highlighted row, to match the source workbook.
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5

# Red, medium "box" border drawn around A6:C6 (left edge on A6, right
# edge on C6, top+bottom running the full width), plus a red font.
$box = $ws.Range("A6:C6")
$box.Font.Color = 255
$box.Font.Name = "新細明體"

$box.Borders.Item(9).Weight = -4138
$box.Borders.Item(9).Color = 255

$a6.Borders.Item(10).LineStyle = 0
$b6.Borders.Item(7).LineStyle = 0
$b6.Borders.Item(10).LineStyle = 0
$c6.Borders.Item(7).LineStyle = 0

# Select the new row, mirroring the saved workbook's cursor position.
$ws.Range("A6:C6").Select()
